# Actualiza la base de datos del Estado de Cuenta:
# se eliminan los periodos de mora anteriores del trabajador EDER LUIS
# MEDRANO SOLANO y del trabajador EDILBERTO ANTONIO LAZARO HIGUITA,
# y se agregan los nuevos periodos/valores, reordenando las filas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 16: EDILBERTO ANTONIO LAZARO HIGUITA - periodo 2210
$ws.Range("C16").Value = "1143342317"
$ws.Range("D16").Value = "EDILBERTO ANTONIO LAZARO HIGUITA"
$ws.Range("E16").Value = "2210"
$ws.Range("F16").Value = 2667

# Fila 17: EDILBERTO ANTONIO LAZARO HIGUITA - periodo 2209
$ws.Range("C17").Value = "1143342317"
$ws.Range("D17").Value = "EDILBERTO ANTONIO LAZARO HIGUITA"
$ws.Range("E17").Value = "2209"
$ws.Range("F17").Value = 40000

# Fila 18: EDILBERTO ANTONIO LAZARO HIGUITA - periodo 2208
$ws.Range("C18").Value = "1143342317"
$ws.Range("D18").Value = "EDILBERTO ANTONIO LAZARO HIGUITA"
$ws.Range("E18").Value = "2208"
$ws.Range("F18").Value = 17333

# Fila 19: EDER LUIS MEDRANO SOLANO - periodo 2202
$ws.Range("C19").Value = "1070806421"
$ws.Range("D19").Value = "EDER LUIS MEDRANO SOLANO"
$ws.Range("E19").Value = "2202"
$ws.Range("F19").Value = 4000
